# Finished updating math library for rigid bodies
# Update the diary log entry for "6 marras" (row 25): extend reading range,
# extend afternoon session time, add two new reflection notes, and bump hours.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25 updates -----------------------------------------------------------

# B25 (Kello / time): afternoon session extended
$ws.Range("B25").Value = "9.15-11.15, 12:15-14.15"

# C25 (Oppimisen sisältö / content): read further into the textbook
$ws.Range("C25").Value = "oppikirjasta 162-192 , Implementing mathematics of rotations, "

# F25 (META) is a newly used column on this row - set it before D25 so the
# shared-string table gets the two new strings in the same order Excel would
# emit them (first-seen-in-document order).
$ws.Range("F25").Value = "Tuhti päivä, osa kaavoista kirjoitettiin itse. Kvaternin kohdalla luovutin."

# D25 (Oppimisen laatu / learning quality) - new cell
$ws.Range("D25").Value = "Aaltoillen, asia oli tuttua mutta uuttakin tuli sopivassa suhteessa. Yritin painottaa etenemistä."

# Match the wrapped-text styling already used by the other notes cells on the row
$ws.Range("D25").WrapText = $true
$ws.Range("F25").WrapText = $true

# G25 (Tunnit / hours): longer day -> more hours logged (H3's SUM(G3:G60) recalcs automatically)
$ws.Range("G25").Value = 4

# Row grew taller to fit the wrapped notes
$ws.Rows.Item(25).RowHeight = 43.5

# View state: scrolled further down and selection moved to the new META cell
$ws.Range("F25").Select()
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 100
